$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text value is purely numeric-looking need to be forced to
# remain text (matching the source inline-string cells) without altering the
# cell style index: apply a temporary "@" (text) number format, write the
# value, then restore the Normal style so the saved style index is unchanged.
function Set-TextValue($range, $text) {
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.Style = "Normal"
}

$ws.Range("D2").Value = '61.595.23'
$ws.Range("E2").Value = '  -2.29%  '
$ws.Range("D3").Value = '2.949.59'
$ws.Range("E3").Value = '  -3.35%  '
$ws.Range("E4").Value = '  +0.07%  '
Set-TextValue $ws.Range("D5") '582.58'
$ws.Range("E5").Value = '  -0.53%  '
Set-TextValue $ws.Range("D6") '141.63'
$ws.Range("E6").Value = '  -6.77%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("E8").Value = '  -2.86%  '
$ws.Range("D9").Value = '2.947.24'
$ws.Range("E9").Value = '  -3.52%  '
$ws.Range("E10").Value = '  -5.89%  '
Set-TextValue $ws.Range("D11") '5.70'
$ws.Range("E11").Value = '  -2.61%  '
$ws.Range("E12").Value = '  +1.99%  '
Set-TextValue $ws.Range("D13") '0.0000224'
$ws.Range("E13").Value = '  -4.17%  '
Set-TextValue $ws.Range("D14") '33.96'
$ws.Range("E14").Value = '  -6.26%  '
$ws.Range("E15").Value = '  +1.40%  '
$ws.Range("D16").Value = '3.440.38'
$ws.Range("E16").Value = '  -3.33%  '
$ws.Range("E17").Value = '  -2.24%  '
$ws.Range("D18").Value = '61.635.96'
$ws.Range("E18").Value = '  -2.21%  '
$ws.Range("D19").Value = '2.947.26'
$ws.Range("E19").Value = '  -3.45%  '
Set-TextValue $ws.Range("D20") '448.23'
$ws.Range("E20").Value = '  -6.50%  '
Set-TextValue $ws.Range("D21") '13.79'
$ws.Range("E21").Value = '  -3.48%  '
Set-TextValue $ws.Range("D22") '0.677'
$ws.Range("E22").Value = '  -4.18%  '
Set-TextValue $ws.Range("D23") '7.25'
$ws.Range("E23").Value = '  -3.36%  '
Set-TextValue $ws.Range("D24") '81.03'
$ws.Range("E24").Value = '  -1.04%  '
$ws.Range("E25").Value = '  -4.42%  '
$ws.Range("E26").Value = '  -10.88%  '
$ws.Range("E27").Value = '  -0.07%  '
Set-TextValue $ws.Range("D28") '9.44'
$ws.Range("E28").Value = '  -10.59%  '
$ws.Range("E29").Value = '  -0.04%  '
$ws.Range("E30").Value = '  -2.17%  '
$ws.Range("E31").Value = '  -7.54%  '
$ws.Range("E32").Value = '  -6.75%  '
Set-TextValue $ws.Range("D33") '27.16'
$ws.Range("E33").Value = '  -1.41%  '
$ws.Range("E34").Value = '  -4.08%  '
Set-TextValue $ws.Range("D35") '1.00'
$ws.Range("E35").Value = '  -5.49%  '
$ws.Range("D36").Value = '0.0₃0768'
$ws.Range("E36").Value = '  -5.99%  '
$ws.Range("E37").Value = '  -4.51%  '
$ws.Range("E38").Value = '  -6.30%  '
Set-TextValue $ws.Range("D39") '49.95'
$ws.Range("E39").Value = '  -0.89%  '
Set-TextValue $ws.Range("D40") '9.08'
$ws.Range("E40").Value = '  -1.73%  '
$ws.Range("E41").Value = '  +3.34%  '
$ws.Range("E42").Value = '  -14.29%  '
Set-TextValue $ws.Range("D43") '387.41'
$ws.Range("E43").Value = '  -9.88%  '
$ws.Range("E44").Value = '  -3.05%  '
$ws.Range("D45").Value = '2.707.50'
$ws.Range("E45").Value = '  -4.47%  '
$ws.Range("E46").Value = '  -9.16%  '
Set-TextValue $ws.Range("D47") '36.58'
$ws.Range("E47").Value = '  -4.25%  '
Set-TextValue $ws.Range("D48") '129.79'
$ws.Range("E48").Value = '  +1.80%  '
$ws.Range("E50").Value = '  -1.94%  '
$ws.Range("E51").Value = '  -2.31%  '
